# Updated cryptos list values (Price and Volume(1h) columns) per upstream scrape refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "60.927.76" }
    @{ Cell = "E2"; Value = "  -2.46%  " }
    @{ Cell = "D3"; Value = "2.422.39" }
    @{ Cell = "E3"; Value = "  -1.36%  " }
    @{ Cell = "D4"; Value = "0.997" }
    @{ Cell = "E4"; Value = "  -0.20%  " }
    @{ Cell = "D5"; Value = "571.42" }
    @{ Cell = "E5"; Value = "  -1.46%  " }
    @{ Cell = "D6"; Value = "139.98" }
    @{ Cell = "E6"; Value = "  -2.64%  " }
    @{ Cell = "E7"; Value = "  +0.20%  " }
    @{ Cell = "E8"; Value = "  -1.25%  " }
    @{ Cell = "D9"; Value = "2.406.83" }
    @{ Cell = "E9"; Value = "  -1.71%  " }
    @{ Cell = "E10"; Value = "  -1.13%  " }
    @{ Cell = "E11"; Value = "  -0.21%  " }
    @{ Cell = "D13"; Value = "0.340" }
    @{ Cell = "E13"; Value = "  -2.20%  " }
    @{ Cell = "D14"; Value = "25.97" }
    @{ Cell = "E14"; Value = "  -1.97%  " }
    @{ Cell = "E15"; Value = "  -4.11%  " }
    @{ Cell = "D16"; Value = "2.829.88" }
    @{ Cell = "E16"; Value = "  -2.43%  " }
    @{ Cell = "D17"; Value = "60.734.43" }
    @{ Cell = "E17"; Value = "  -2.47%  " }
    @{ Cell = "D18"; Value = "2.422.76" }
    @{ Cell = "E18"; Value = "  -1.16%  " }
    @{ Cell = "D19"; Value = "7.56" }
    @{ Cell = "E19"; Value = "  +5.77%  " }
    @{ Cell = "D20"; Value = "10.70" }
    @{ Cell = "E20"; Value = "  -1.70%  " }
    @{ Cell = "D21"; Value = "322.78" }
    @{ Cell = "E21"; Value = "  -1.89%  " }
    @{ Cell = "E22"; Value = "  -2.00%  " }
    @{ Cell = "E23"; Value = "  +0.82%  " }
    @{ Cell = "E24"; Value = "  +0.11%  " }
    @{ Cell = "E25"; Value = "  -5.24%  " }
    @{ Cell = "D26"; Value = "64.86" }
    @{ Cell = "E26"; Value = "  -1.33%  " }
    @{ Cell = "D27"; Value = "586.17" }
    @{ Cell = "E27"; Value = "  -0.39%  " }
    @{ Cell = "D28"; Value = "8.44" }
    @{ Cell = "E28"; Value = "  -9.52%  " }
    @{ Cell = "D29"; Value = "2.526.50" }
    @{ Cell = "E29"; Value = "  -1.93%  " }
    @{ Cell = "D30"; Value = "0.0₃0924" }
    @{ Cell = "E30"; Value = "  -4.60%  " }
    @{ Cell = "D31"; Value = "7.93" }
    @{ Cell = "E31"; Value = "  -1.70%  " }
    @{ Cell = "D32"; Value = "1.36" }
    @{ Cell = "E32"; Value = "  -5.65%  " }
    @{ Cell = "E33"; Value = "  -3.05%  " }
    @{ Cell = "E34"; Value = "  -1.60%  " }
    @{ Cell = "E35"; Value = "  -0.04%  " }
    @{ Cell = "D36"; Value = "4.64" }
    @{ Cell = "E36"; Value = "  -6.25%  " }
    @{ Cell = "E37"; Value = "  -3.40%  " }
    @{ Cell = "D38"; Value = "151.72" }
    @{ Cell = "E38"; Value = "  -2.04%  " }
    @{ Cell = "E39"; Value = "  -3.19%  " }
    @{ Cell = "D40"; Value = "18.28" }
    @{ Cell = "E40"; Value = "  -0.96%  " }
    @{ Cell = "E41"; Value = "  -3.77%  " }
    @{ Cell = "E42"; Value = "  +0.11%  " }
    @{ Cell = "E43"; Value = "  -2.95%  " }
    @{ Cell = "D44"; Value = "41.19" }
    @{ Cell = "D45"; Value = "2.36" }
    @{ Cell = "E45"; Value = "  -4.56%  " }
    @{ Cell = "D46"; Value = "0.0₆0287" }
    @{ Cell = "E46"; Value = "  +9.62%  " }
    @{ Cell = "D47"; Value = "141.18" }
    @{ Cell = "E47"; Value = "  -1.49%  " }
    @{ Cell = "E48"; Value = "  -4.14%  " }
    @{ Cell = "D49"; Value = "0.591" }
    @{ Cell = "E49"; Value = "  -3.15%  " }
    @{ Cell = "D50"; Value = "19.60" }
    @{ Cell = "E50"; Value = "  -1.68%  " }
    @{ Cell = "D51"; Value = "0.0504" }
    @{ Cell = "E51"; Value = "  -3.92%  " }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
    $rng.Style = "Normal"
}
